$d = $word.ActiveDocument

# Append a new paragraph at the end of the document (after the
# "... Flop." entry) containing the "30 de agosto" bitácora entry.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 12
$newRange.LanguageID = "es-ES"

# Insert the two runs ("30 de agosto: " and the rest of the sentence) as
# two distinct tracked-change insertions so they remain separate <w:r>
# elements (with identical run formatting) instead of being coalesced
# into a single run, then accept each insertion individually so the
# run-level formatting (rFonts/sz/szCs/lang) is fully materialized.
$d.TrackRevisions = $true
$newRange.InsertAfter("30 de agosto: ")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

$secondRange = $d.Paragraphs.Last.Range
$secondRange.Collapse(0)
$secondRange.Font.Name = "Times New Roman"
$secondRange.Font.Size = 12
$secondRange.LanguageID = "es-ES"

$d.TrackRevisions = $true
$secondRange.InsertAfter("Después de investigar sobre los integrados menos comunes, se llega a la conclusión de que es posible simplificar aun más el circuito del acumulador si se utiliza compuertas como XOR y XNOR, gracias a la simplificación por medio de propiedades de algebra booleana. ")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()
